# The presentation ships two theme parts:
#   ppt/theme/theme1.xml  ("Office Theme" colours) - used by the notes master
#   ppt/theme/theme2.xml  ("Integral" colours)      - used by the slide master / design
#
# The authored edit swaps the two themes' contents: the design that drives the
# slides (currently "Integral") takes on the old "Office Theme" palette, and
# vice-versa. The live design's theme colours are reachable through
# Slide.ThemeColorScheme.Colors(i).RGB, matching the clrScheme child order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

function ComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

# Target palette: the "Office Theme" colours that used to live in theme1.xml,
# now applied (after the swap) to the design backing the slides.
$newColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p   = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $newColors.Count; $i++) {
    $tcs.Colors($i).RGB = ComRGB($newColors[$i - 1])
}
